$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F. Excel shifts the existing column F
# (the "genus/species/family/NA" unit labels) to column G, carrying its
# formatting along, and the new column F inherits the same per-row styles.
$ws.Columns("F:F").Insert()

# New "Feb" header in the freshly inserted column.
$ws.Range("F2").Value = "Feb"

# New Feb data values for each row.
$ws.Range("F3").Value = 136354
$ws.Range("F4").Value = 10750
$ws.Range("F5").Value = 13530
$ws.Range("F6").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("F8").Value = 2115
$ws.Range("F9").Value = 6919

$ws.Range("F11").Formula = "=F8/F9"

$ws.Range("F13").Value = 17628
$ws.Range("F14").Value = "NA"
$ws.Range("F15").Value = 16747
$ws.Range("F16").Value = 132

# Match the saved selection/active range.
$ws.Range("A2:G16").Select()
